$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.786.55"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "2.701.19"
$ws.Range("E3").Value = "  +2.93%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +5.86%  "
$ws.Range("E10").Value = "  +4.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.402"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000203"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.33%  "
$ws.Range("D15").Value = "3.189.32"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").Value = "65.703.88"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "2.693.71"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "359.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.65%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.43%  "
$ws.Range("E25").Value = "  +13.07%  "
$ws.Range("E26").Value = "  +2.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.35%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.171"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.51%  "
$ws.Range("E30").Value = "  +5.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "545.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.42%  "
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +3.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.99%  "
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.434"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "172.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.11%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.00%  "
$ws.Range("E44").Value = "  +3.50%  "
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.99%  "
$ws.Range("E48").Value = "  +5.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.655"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0993"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.66%  "

Write-Output "Applied 87 cell updates"